$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 1 (Queue Size vs deadlocks, rows 3-9): updated trial counts.
#     Column D holds =((Bn+Cn)/2)*20 and recomputes automatically. ---
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 3

$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 2

$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 1

$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 1

$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 0

$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0

# --- Table 2 (Consumers vs deadlocks, rows 34-40): trial counts filled in
#     for the first time (previously blank, formulas evaluated to 0). ---
$ws.Range("B34").Value = 0
$ws.Range("C34").Value = 0

$ws.Range("B35").Value = 0
$ws.Range("C35").Value = 0

$ws.Range("B36").Value = 0
$ws.Range("C36").Value = 1

$ws.Range("B37").Value = 2
$ws.Range("C37").Value = 2

$ws.Range("B38").Value = 3
$ws.Range("C38").Value = 4

$ws.Range("B39").Value = 5
$ws.Range("C39").Value = 5

$ws.Range("B40").Value = 5
$ws.Range("C40").Value = 5

# --- Move the view: scroll down to the second table and select F41 ---
$null = $ws.Range("F41").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
